$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update run_time (C) and max_er/iter columns (E:Y) for rows 2-11
# with newly simulated results for the other-language run

$ws.Range("C2").Value = 0.3859937191009521
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 166.7008502018489
$rowVals[0,1] = 0.005695122577285056
$rowVals[0,2] = 0.004889840092523848
$rowVals[0,3] = 0.004487606760621414
$rowVals[0,4] = 0.00440117539087203
$rowVals[0,5] = 0.004137959625192721
$rowVals[0,6] = 0.004015817776864572
$rowVals[0,7] = 0.003850441125197786
$rowVals[0,8] = 0.003683872805420258
$rowVals[0,9] = 0.003683872805420258
$rowVals[0,10] = 0.003683265085562939
$rowVals[0,11] = 0.00360966234782559
$rowVals[0,12] = 0.00351014634416858
$rowVals[0,13] = 0.003493996494529795
$rowVals[0,14] = 0.003448206863847684
$rowVals[0,15] = 0.003346620104944876
$rowVals[0,16] = 0.003331217051239995
$rowVals[0,17] = 0.003315106628683085
$rowVals[0,18] = 0.003258553129899101
$rowVals[0,19] = 0.003258553129899101
$rowVals[0,20] = 0.003249529243700757
$ws.Range("E2:Y2").Value = $rowVals

$ws.Range("C3").Value = 0.4240038394927979
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 168.4844816457371
$rowVals[0,1] = 0.005547066448613553
$rowVals[0,2] = 0.004815790443408508
$rowVals[0,3] = 0.004359210186529429
$rowVals[0,4] = 0.004044440135451783
$rowVals[0,5] = 0.003930571162451782
$rowVals[0,6] = 0.003827015977974235
$rowVals[0,7] = 0.003827015977974235
$rowVals[0,8] = 0.00380315303388296
$rowVals[0,9] = 0.00380315303388296
$rowVals[0,10] = 0.003586118906002835
$rowVals[0,11] = 0.003586118906002835
$rowVals[0,12] = 0.003514412508183452
$rowVals[0,13] = 0.003514412508183452
$rowVals[0,14] = 0.003377760313186951
$rowVals[0,15] = 0.003373852811136243
$rowVals[0,16] = 0.003366210280398758
$rowVals[0,17] = 0.00335320082596479
$rowVals[0,18] = 0.003330886865575474
$rowVals[0,19] = 0.003297363480766175
$rowVals[0,20] = 0.003284297887831132
$ws.Range("E3:Y3").Value = $rowVals

$ws.Range("C4").Value = 0.3389937877655029
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 172.9706161829999
$rowVals[0,1] = 0.005642025233686275
$rowVals[0,2] = 0.004697614393623904
$rowVals[0,3] = 0.004222517796814512
$rowVals[0,4] = 0.004207422172915288
$rowVals[0,5] = 0.003942887034631921
$rowVals[0,6] = 0.003942887034631921
$rowVals[0,7] = 0.003942887034631921
$rowVals[0,8] = 0.003795615851558623
$rowVals[0,9] = 0.003795615851558623
$rowVals[0,10] = 0.003752079240196846
$rowVals[0,11] = 0.003665740883259158
$rowVals[0,12] = 0.003632893779542026
$rowVals[0,13] = 0.003632893779542026
$rowVals[0,14] = 0.003564690947618674
$rowVals[0,15] = 0.003493115820322583
$rowVals[0,16] = 0.00344813076534484
$rowVals[0,17] = 0.003414444108895115
$rowVals[0,18] = 0.003414444108895115
$rowVals[0,19] = 0.003388558051688076
$rowVals[0,20] = 0.003371746904152045
$ws.Range("E4:Y4").Value = $rowVals

$ws.Range("C5").Value = 0.3330380916595459
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 164.144334356024
$rowVals[0,1] = 0.005679871177659367
$rowVals[0,2] = 0.004647144724443379
$rowVals[0,3] = 0.004350897813186371
$rowVals[0,4] = 0.004295348244944137
$rowVals[0,5] = 0.004060491610429554
$rowVals[0,6] = 0.003849562872211274
$rowVals[0,7] = 0.003837585557495068
$rowVals[0,8] = 0.003755720966883497
$rowVals[0,9] = 0.003713554925785292
$rowVals[0,10] = 0.003677229587532313
$rowVals[0,11] = 0.00349717933042899
$rowVals[0,12] = 0.003409881610060627
$rowVals[0,13] = 0.003398182646672929
$rowVals[0,14] = 0.003398182646672929
$rowVals[0,15] = 0.003354594513393376
$rowVals[0,16] = 0.0033336010553651
$rowVals[0,17] = 0.003282845760687374
$rowVals[0,18] = 0.003237313432479307
$rowVals[0,19] = 0.00319969462682308
$rowVals[0,20] = 0.00319969462682308
$ws.Range("E5:Y5").Value = $rowVals

$ws.Range("C6").Value = 0.3959712982177734
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 169.453223336699
$rowVals[0,1] = 0.005643394092958738
$rowVals[0,2] = 0.004622255019624994
$rowVals[0,3] = 0.004316892915374845
$rowVals[0,4] = 0.004241624078074937
$rowVals[0,5] = 0.003883221527374979
$rowVals[0,6] = 0.003883221527374979
$rowVals[0,7] = 0.003746426981202571
$rowVals[0,8] = 0.00364105728197811
$rowVals[0,9] = 0.003612257140037286
$rowVals[0,10] = 0.003557360060473909
$rowVals[0,11] = 0.003522268380618743
$rowVals[0,12] = 0.003498995550691783
$rowVals[0,13] = 0.003498995550691783
$rowVals[0,14] = 0.003447557648242333
$rowVals[0,15] = 0.003393724237416001
$rowVals[0,16] = 0.003360912242698961
$rowVals[0,17] = 0.003331715421806725
$rowVals[0,18] = 0.003317364178354565
$rowVals[0,19] = 0.00330318174145612
$rowVals[0,20] = 0.00330318174145612
$ws.Range("E6:Y6").Value = $rowVals

$ws.Range("C7").Value = 0.3299975395202637
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 171.1355819736091
$rowVals[0,1] = 0.005387905113368413
$rowVals[0,2] = 0.004500859285083696
$rowVals[0,3] = 0.004332660443315584
$rowVals[0,4] = 0.004114736593981621
$rowVals[0,5] = 0.004114736593981621
$rowVals[0,6] = 0.003959596441489703
$rowVals[0,7] = 0.003867890377220679
$rowVals[0,8] = 0.003840430743196965
$rowVals[0,9] = 0.003723555608021183
$rowVals[0,10] = 0.003592625512288332
$rowVals[0,11] = 0.003561733560511189
$rowVals[0,12] = 0.003552385442953093
$rowVals[0,13] = 0.003529394580070009
$rowVals[0,14] = 0.003529394580070009
$rowVals[0,15] = 0.003417676049170355
$rowVals[0,16] = 0.003417676049170355
$rowVals[0,17] = 0.003388864710747641
$rowVals[0,18] = 0.003353699098812749
$rowVals[0,19] = 0.00334191513817284
$rowVals[0,20] = 0.003335976256795498
$ws.Range("E7:Y7").Value = $rowVals

$ws.Range("C8").Value = 0.4490087032318115
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 170.0916393539883
$rowVals[0,1] = 0.005695122577285056
$rowVals[0,2] = 0.004830976227410687
$rowVals[0,3] = 0.004381125642746454
$rowVals[0,4] = 0.004105635070414723
$rowVals[0,5] = 0.004091818112775553
$rowVals[0,6] = 0.003874774837468864
$rowVals[0,7] = 0.003785155475214207
$rowVals[0,8] = 0.003785155475214207
$rowVals[0,9] = 0.003662268242007509
$rowVals[0,10] = 0.003662268242007509
$rowVals[0,11] = 0.003651657548090693
$rowVals[0,12] = 0.003558577899002557
$rowVals[0,13] = 0.00355299168544906
$rowVals[0,14] = 0.003492726961722991
$rowVals[0,15] = 0.003418610263837564
$rowVals[0,16] = 0.003394464542446808
$rowVals[0,17] = 0.003358325462315661
$rowVals[0,18] = 0.003330151004204018
$rowVals[0,19] = 0.003327475252368301
$rowVals[0,20] = 0.003315626498128426
$ws.Range("E8:Y8").Value = $rowVals

$ws.Range("C9").Value = 0.4629929065704346
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 172.2880020228422
$rowVals[0,1] = 0.005695122577285056
$rowVals[0,2] = 0.004942781660371025
$rowVals[0,3] = 0.004421185863639588
$rowVals[0,4] = 0.004310221364566313
$rowVals[0,5] = 0.004098535484569515
$rowVals[0,6] = 0.004098535484569515
$rowVals[0,7] = 0.003824420124755134
$rowVals[0,8] = 0.003794128503775557
$rowVals[0,9] = 0.003794128503775557
$rowVals[0,10] = 0.003625737247789144
$rowVals[0,11] = 0.003594895334227209
$rowVals[0,12] = 0.003508713540794696
$rowVals[0,13] = 0.003508713540794696
$rowVals[0,14] = 0.003507537585664125
$rowVals[0,15] = 0.003475111910277325
$rowVals[0,16] = 0.003404520808540894
$rowVals[0,17] = 0.003404520808540894
$rowVals[0,18] = 0.003395322483326689
$rowVals[0,19] = 0.003368031084317135
$rowVals[0,20] = 0.003358440585240588
$ws.Range("E9:Y9").Value = $rowVals

$ws.Range("C10").Value = 0.4579970836639404
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 165.5295241311142
$rowVals[0,1] = 0.005695122577285056
$rowVals[0,2] = 0.00479843212934127
$rowVals[0,3] = 0.004659000095937505
$rowVals[0,4] = 0.004216417665214705
$rowVals[0,5] = 0.00387386193980597
$rowVals[0,6] = 0.00387386193980597
$rowVals[0,7] = 0.003839103593039638
$rowVals[0,8] = 0.00364647870247139
$rowVals[0,9] = 0.003424646009772371
$rowVals[0,10] = 0.003424646009772371
$rowVals[0,11] = 0.003424646009772371
$rowVals[0,12] = 0.003424646009772371
$rowVals[0,13] = 0.003411907271017616
$rowVals[0,14] = 0.003373152794837272
$rowVals[0,15] = 0.003360712189673842
$rowVals[0,16] = 0.00325971926022195
$rowVals[0,17] = 0.00325971926022195
$rowVals[0,18] = 0.00325971926022195
$rowVals[0,19] = 0.003255445824384349
$rowVals[0,20] = 0.003226696376824837
$ws.Range("E10:Y10").Value = $rowVals

$ws.Range("C11").Value = 0.478001594543457
$rowVals = New-Object 'object[,]' 1,21
$rowVals[0,0] = 168.1139107697763
$rowVals[0,1] = 0.005695122577285056
$rowVals[0,2] = 0.004900069624440375
$rowVals[0,3] = 0.004631824322849579
$rowVals[0,4] = 0.00437560012312078
$rowVals[0,5] = 0.004136599001335569
$rowVals[0,6] = 0.004124312046736903
$rowVals[0,7] = 0.003728057260053844
$rowVals[0,8] = 0.003700556898193613
$rowVals[0,9] = 0.003700556898193613
$rowVals[0,10] = 0.003552893236305676
$rowVals[0,11] = 0.003552893236305676
$rowVals[0,12] = 0.003552893236305676
$rowVals[0,13] = 0.00341246171402593
$rowVals[0,14] = 0.00341246171402593
$rowVals[0,15] = 0.00340877921095578
$rowVals[0,16] = 0.003370549890403222
$rowVals[0,17] = 0.003298251570178894
$rowVals[0,18] = 0.003298251570178894
$rowVals[0,19] = 0.003286693239026301
$rowVals[0,20] = 0.003277074284011234
$ws.Range("E11:Y11").Value = $rowVals
